$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Rename the diff-comparison column headers:
#   "<Field>_old" -> "<Field>_FV2310"
#   "<Field>_new" -> "<Field>_FV2404"
$fields = @("Segmentname", "Segmentgruppe", "Segment", "Datenelement", "Segment ID", "Code", "Qualifier", "Beschreibung", "Bedingungsausdruck", "Bedingung")

for ($i = 0; $i -lt $fields.Length; $i++) {
    # Columns A..J (1..10) hold the "_old" headers -> rename to "_FV2310"
    $ws.Cells.Item(1, $i + 1).Value = $fields[$i] + "_FV2310"
    # Columns L..U (12..21) hold the "_new" headers -> rename to "_FV2404"
    $ws.Cells.Item(1, $i + 12).Value = $fields[$i] + "_FV2404"
}
# Column K (11) keeps its header "diff" - unchanged.

# Turn the used range into an Excel Table ("Table1") with an autofilter,
# matching the header row that was just renamed.
$rng = $ws.Range("A1:U92")
$tbl = $ws.ListObjects.Add(1, $rng, 0, 1)
$tbl.Name = "Table1"
$tbl.TableStyle = ""

# Freeze the header row (split after row 1): select the first cell below the
# header and freeze so the top pane (row 1) stays visible while scrolling.
$ws.Range("A2").Select()
$ws.Application.ActiveWindow.FreezePanes = $true
